$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Programa") and add the "Facultad" header
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "Facultad"

# Update selection to match the recorded state (entire column E selected)
$ws.Range("E1:E1048576").Select()
